$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Permuted species-record block (rows 2-12, columns A,B,D,E,F,G,H,Q,R,S)
# Row-local columns (C,I,P,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AT,AW,AX,AY) are unchanged.
# Row 2
$ws.Cells.Item(2, 1).Value = 111363029
$ws.Cells.Item(2, 2).Value = 76918
$ws.Cells.Item(2, 4).Value = 'NT'
$ws.Cells.Item(2, 5).Value = 6437
$ws.Cells.Item(2, 6).Value = 'Blanksvart spiklav'
$ws.Cells.Item(2, 7).Value = 'Calicium denigratum'
$ws.Cells.Item(2, 8).Value = '(Vain.) Tibell'
$ws.Cells.Item(2, 17).Value = 593312.9580448985
$ws.Cells.Item(2, 18).Value = 6987010.291132212
$ws.Cells.Item(2, 19).Value = 10

# Row 3
$ws.Cells.Item(3, 1).Value = 111363023
$ws.Cells.Item(3, 2).Value = 76918
$ws.Cells.Item(3, 4).Value = 'NT'
$ws.Cells.Item(3, 5).Value = 6437
$ws.Cells.Item(3, 6).Value = 'Blanksvart spiklav'
$ws.Cells.Item(3, 7).Value = 'Calicium denigratum'
$ws.Cells.Item(3, 8).Value = '(Vain.) Tibell'
$ws.Cells.Item(3, 17).Value = 593269.3631576585
$ws.Cells.Item(3, 18).Value = 6987149.513888635
$ws.Cells.Item(3, 19).Value = 10

# Row 4
$ws.Cells.Item(4, 1).Value = 111363030
$ws.Cells.Item(4, 2).Value = 77268
$ws.Cells.Item(4, 4).Value = 'NT'
$ws.Cells.Item(4, 5).Value = 228912
$ws.Cells.Item(4, 6).Value = 'Mörk kolflarnlav'
$ws.Cells.Item(4, 7).Value = 'Carbonicola myrmecina'
$ws.Cells.Item(4, 8).Value = '(Ach.) Bendiksby & Timdal'
$ws.Cells.Item(4, 17).Value = 593355.1995546351
$ws.Cells.Item(4, 18).Value = 6987156.520171621
$ws.Cells.Item(4, 19).Value = 25

# Row 5
$ws.Cells.Item(5, 1).Value = 111363024
$ws.Cells.Item(5, 2).Value = 77268
$ws.Cells.Item(5, 4).Value = 'NT'
$ws.Cells.Item(5, 5).Value = 228912
$ws.Cells.Item(5, 6).Value = 'Mörk kolflarnlav'
$ws.Cells.Item(5, 7).Value = 'Carbonicola myrmecina'
$ws.Cells.Item(5, 8).Value = '(Ach.) Bendiksby & Timdal'
$ws.Cells.Item(5, 17).Value = 593291.0260186956
$ws.Cells.Item(5, 18).Value = 6987171.95495991
$ws.Cells.Item(5, 19).Value = 10

# Row 6
$ws.Cells.Item(6, 1).Value = 111363031
$ws.Cells.Item(6, 2).Value = 76918
$ws.Cells.Item(6, 4).Value = 'NT'
$ws.Cells.Item(6, 5).Value = 6437
$ws.Cells.Item(6, 6).Value = 'Blanksvart spiklav'
$ws.Cells.Item(6, 7).Value = 'Calicium denigratum'
$ws.Cells.Item(6, 8).Value = '(Vain.) Tibell'
$ws.Cells.Item(6, 17).Value = 593417.4633552339
$ws.Cells.Item(6, 18).Value = 6986985.556671137
$ws.Cells.Item(6, 19).Value = 10

# Row 7
$ws.Cells.Item(7, 1).Value = 111363020
$ws.Cells.Item(7, 2).Value = 78107
$ws.Cells.Item(7, 4).Value = 'NT'
$ws.Cells.Item(7, 5).Value = 6453
$ws.Cells.Item(7, 6).Value = 'Vedskivlav'
$ws.Cells.Item(7, 7).Value = 'Hertelidea botryosa'
$ws.Cells.Item(7, 8).Value = '(Fr.) Printzen & Kantvilas'
$ws.Cells.Item(7, 17).Value = 593324.7367794912
$ws.Cells.Item(7, 18).Value = 6987171.102828567
$ws.Cells.Item(7, 19).Value = 10

# Row 8
$ws.Cells.Item(8, 1).Value = 111363021
$ws.Cells.Item(8, 2).Value = 89330
$ws.Cells.Item(8, 4).Value = 'NT'
$ws.Cells.Item(8, 5).Value = 3242
$ws.Cells.Item(8, 6).Value = 'Vitplätt'
$ws.Cells.Item(8, 7).Value = 'Chaetodermella luna'
$ws.Cells.Item(8, 8).Value = '(Romell ex D.P.Rogers & H.S.Jacks.) Rauschert'
$ws.Cells.Item(8, 17).Value = 593278.356042281
$ws.Cells.Item(8, 18).Value = 6987153.408284122
$ws.Cells.Item(8, 19).Value = 10

# Row 9
$ws.Cells.Item(9, 1).Value = 111363022
$ws.Cells.Item(9, 2).Value = 77186
$ws.Cells.Item(9, 4).Value = 'NT'
$ws.Cells.Item(9, 5).Value = 353
$ws.Cells.Item(9, 6).Value = 'Dvärgbägarlav'
$ws.Cells.Item(9, 7).Value = 'Cladonia parasitica'
$ws.Cells.Item(9, 8).Value = '(Hoffm.) Hoffm.'
$ws.Cells.Item(9, 17).Value = 593324.9051589288
$ws.Cells.Item(9, 18).Value = 6987181.108611984
$ws.Cells.Item(9, 19).Value = 10

# Row 10
$ws.Cells.Item(10, 1).Value = 111363025
$ws.Cells.Item(10, 2).Value = 89646
$ws.Cells.Item(10, 4).Value = 'VU'
$ws.Cells.Item(10, 5).Value = 65
$ws.Cells.Item(10, 6).Value = 'Fläckporing'
$ws.Cells.Item(10, 7).Value = 'Anthoporia albobrunnea'
$ws.Cells.Item(10, 8).Value = '(Romell) Karasiński & Niemelä'
$ws.Cells.Item(10, 17).Value = 593292.3890792141
$ws.Cells.Item(10, 18).Value = 6987203.815111163
$ws.Cells.Item(10, 19).Value = 10

# Row 11
$ws.Cells.Item(11, 1).Value = 111363028
$ws.Cells.Item(11, 2).Value = 77186
$ws.Cells.Item(11, 4).Value = 'NT'
$ws.Cells.Item(11, 5).Value = 353
$ws.Cells.Item(11, 6).Value = 'Dvärgbägarlav'
$ws.Cells.Item(11, 7).Value = 'Cladonia parasitica'
$ws.Cells.Item(11, 8).Value = '(Hoffm.) Hoffm.'
$ws.Cells.Item(11, 17).Value = 593324.0129203054
$ws.Cells.Item(11, 18).Value = 6987101.07452714
$ws.Cells.Item(11, 19).Value = 10

# Row 12
$ws.Cells.Item(12, 1).Value = 111363026
$ws.Cells.Item(12, 2).Value = 90854
$ws.Cells.Item(12, 4).Value = 'NT'
$ws.Cells.Item(12, 5).Value = 2079
$ws.Cells.Item(12, 6).Value = 'Nordtagging'
$ws.Cells.Item(12, 7).Value = 'Odonticium romellii'
$ws.Cells.Item(12, 8).Value = '(S.Lundell) Parmasto'
$ws.Cells.Item(12, 17).Value = 593292.3890792141
$ws.Cells.Item(12, 18).Value = 6987203.815111163
$ws.Cells.Item(12, 19).Value = 10

